# "retour arriere sur les modif écrans" — revert the screen-list edits on
# sheet LIST: restore A2:A4 to their earlier (shorter) values, clear A5:A8,
# move the row style that used to sit on A8 up onto A4 (and vice versa),
# and put the active selection back on A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

# --- values -----------------------------------------------------------
$ws.Range("A2").Value = "AD.SEC.001.FON.01"
$ws.Range("A3").Value = "MP.CPT"
$ws.Range("A4").Value = "AD.SEC.014.FON.01"

$ws.Range("A5").ClearContents()
$ws.Range("A6").ClearContents()
$ws.Range("A7").ClearContents()
$ws.Range("A8").ClearContents()

# --- formatting ---------------------------------------------------------
# A4 loses the "@" text format it had (back to the sheet's default/General)
$ws.Range("A4").Style = "Normal"
# A8 gains the "@" text format that the other list rows use
$ws.Range("A8").NumberFormat = "@"

# --- selection ------------------------------------------------------------
$ws.Activate()
$ws.Range("A8").Select()
